$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) New BOM line item in row 41: Hole Saw 5/8ths ---
# (the vendor URL is interned into the shared-string table first, matching
#  the original author's edit order)
$holeSawUrl = "https://www.grainger.com/product/MORSE-Hole-Saw-5-8-in-Saw-Dia-54HP10"
$ws.Range("G42").Value = $holeSawUrl

$ws.Range("A41").Value = "Hole Saw 5/8ths"
$ws.Range("B41").Value = "MORSE Hole Saw: 5/8 in Saw Dia., 5/6 Teeth per Inch, 1 7/8 in Max. Cutting Dp, 1/2""-20 Thread Size"
$ws.Range("C41").Value = 1
$ws.Range("D41").Value = 23.73
$ws.Range("E41").Value = 23.73

$ws.Range("A40:D40").Copy() | Out-Null
$ws.Range("A41:D41").PasteSpecial(-4122) | Out-Null

# --- 2) Move the grand-total row down: old row 42 (Total) becomes new row 43 ---
$ws.Range("D43").Value = "Total"
$ws.Range("D42").Copy() | Out-Null
$ws.Range("D43").PasteSpecial(-4122) | Out-Null

$ws.Range("E43").Formula = "=SUM(E3:E42)"
$ws.Range("E42").Copy() | Out-Null
$ws.Range("E43").PasteSpecial(-4122) | Out-Null

# --- 3) Old row 42 is cleared of its former Total content ... ---
$ws.Range("D42:F42").ClearContents()

# --- ... and repurposed to hold a Notes hyperlink to the new part's vendor page ---
$ws.Hyperlinks.Add($ws.Range("G42"), $holeSawUrl) | Out-Null
$ws.Range("G40").Copy() | Out-Null
$ws.Range("G42").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
$ws.Range("B43").Select()
